# Weekly fairs handling change: reorder tied-count goods labels in column A.
# Counts in column B stay the same for each row; only the text label in
# column A is being rearranged among rows that share the same count value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "серебреный товар"
$ws.Range("A18").Value = "деревенский товар"

$ws.Range("A25").Value = "набойчатый товар"
$ws.Range("A26").Value = "нужный товар"

$ws.Range("A28").Value = "недорогой товар"
$ws.Range("A29").Value = "внутренний товар"
$ws.Range("A30").Value = "суровский товар"
$ws.Range("A31").Value = "медный товар"
$ws.Range("A32").Value = "питейный припасы"

$ws.Range("A33").Value = "оловянный товар"
$ws.Range("A34").Value = "привозный товар"

$ws.Range("A37").Value = "галантерейный товар"
$ws.Range("A38").Value = "заморский товар"

$ws.Range("A39").Value = "харчевой припасы"
$ws.Range("A40").Value = "меховой товар"
$ws.Range("A41").Value = "домовый товар"
